$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns (zh-cn / de-de) for the handback entry ---
$wsOverview.Range("E2").Value() = $newStatus
$wsOverview.Range("F2").Value() = $newStatus

# --- zh-cn sheet: Status, Latest Handback DateTime, Error Detail ---
$wsZhCn.Range("C2").Value() = $newStatus
$wsZhCn.Range("K2").Value() = "2016-10-18 11:31:54"
$wsZhCn.Range("P2").Value() = ""

# --- de-de sheet: Status, Latest Handback DateTime, Error Detail ---
$wsDeDe.Range("C2").Value() = $newStatus
$wsDeDe.Range("K2").Value() = "2016-10-18 11:32:26"
$wsDeDe.Range("P2").Value() = ""

# --- Column width adjustments (Status column widened, Error Detail column narrowed) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333332

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333332
